$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above current row 5 (shifts existing rows 5:29 down to 7:31,
# inheriting the row-4 formatting the same way Excel does natively).
$ws.Range("A5:A6").EntireRow.Insert()

# New weekly observation dated 2022-05-04 ("Especial" and "Primera" quality rows)
$ws.Range("A5").Value = 8
$ws.Range("B5").Value = "Terminal La Palmera de La Serena"
$ws.Range("C5").Value = "Coquimbo"
$ws.Range("D5").Value = 44685
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100107
$ws.Range("H5").Value = "Otros"
$ws.Range("I5").Value = 100107011
$ws.Range("J5").Value = "Tuna"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Especial"
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 19000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 19500
$ws.Range("Q5").Value = "`$/caja 18 kilos"
$ws.Range("R5").Value = "Provincia de Limarí"
$ws.Range("S5").Value = 1083
$ws.Range("T5").Value = 18

$ws.Range("A6").Value = 8
$ws.Range("B6").Value = "Terminal La Palmera de La Serena"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("D6").Value = 44685
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100107
$ws.Range("H6").Value = "Otros"
$ws.Range("I6").Value = 100107011
$ws.Range("J6").Value = "Tuna"
$ws.Range("K6").Value = "Sin especificar"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 160
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 15500
$ws.Range("Q6").Value = "`$/caja 18 kilos"
$ws.Range("R6").Value = "Provincia de Limarí"
$ws.Range("S6").Value = 861
$ws.Range("T6").Value = 18
